$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "function complete quarantine person" - fill in the remaining
# identification / insurance / arrival-date details for row 2.
$ws.Range("G2").Value = 12222    # CMND/CCCD
$ws.Range("H2").Value = 1        # Ma bao hiem
$ws.Range("M2").Value = 44545    # Ngay den (date, keeps existing date format)

# Reflect where the user ended up looking / the view state at save time.
$win = $excel.ActiveWindow
$win.Zoom = 100
$ws.Range("H2").Select()
